# Fruta / hortaliza, semanal
# Insert two new weekly observation rows into the dataset.
# This shifts all existing rows from 401..507 down to 403..509
# and extends the used range to A1:T509.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 401 (existing data moves down by 2 rows)
$ws.Rows.Item(401).Resize(2).EntireRow.Insert()

# --- Populate new row 401 ---
$ws.Cells.Item(401, 1).Value = 6
$ws.Cells.Item(401, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(401, 3).Value = "Metropolitana"
$ws.Cells.Item(401, 4).Value = 44932
$ws.Cells.Item(401, 5).Value = 13
$ws.Cells.Item(401, 6).Value = "Fruta"
$ws.Cells.Item(401, 7).Value = 100101
$ws.Cells.Item(401, 8).Value = "Berries"
$ws.Cells.Item(401, 9).Value = 100101001
$ws.Cells.Item(401, 10).Value = "Arándano (blue)"
$ws.Cells.Item(401, 11).Value = "Sin especificar"
$ws.Cells.Item(401, 12).Value = "Especial"
$ws.Cells.Item(401, 13).Value = 200
$ws.Cells.Item(401, 14).Value = 3000
$ws.Cells.Item(401, 15).Value = 3000
$ws.Cells.Item(401, 16).Value = 3000
$ws.Cells.Item(401, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(401, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(401, 19).Value = 1500
$ws.Cells.Item(401, 20).Value = 2

# --- Populate new row 402 ---
$ws.Cells.Item(402, 1).Value = 6
$ws.Cells.Item(402, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(402, 3).Value = "Metropolitana"
$ws.Cells.Item(402, 4).Value = 44932
$ws.Cells.Item(402, 5).Value = 13
$ws.Cells.Item(402, 6).Value = "Fruta"
$ws.Cells.Item(402, 7).Value = 100101
$ws.Cells.Item(402, 8).Value = "Berries"
$ws.Cells.Item(402, 9).Value = 100101001
$ws.Cells.Item(402, 10).Value = "Arándano (blue)"
$ws.Cells.Item(402, 11).Value = "Sin especificar"
$ws.Cells.Item(402, 12).Value = "Especial"
$ws.Cells.Item(402, 13).Value = 300
$ws.Cells.Item(402, 14).Value = 3000
$ws.Cells.Item(402, 15).Value = 3000
$ws.Cells.Item(402, 16).Value = 3000
$ws.Cells.Item(402, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(402, 18).Value = "Región del Maule"
$ws.Cells.Item(402, 19).Value = 1500
$ws.Cells.Item(402, 20).Value = 2
